$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 63: 3/19/2010, 2 hours, OMAP UART documentation ---
$ws.Range("A62").Copy() | Out-Null
$ws.Range("A63").PasteSpecial(-4122) | Out-Null
$ws.Range("A63").Value = "03/19/2010"
$ws.Range("B63").Value = 2
$ws.Range("C63").Value = "OMAP UART documentation"

# --- Row 64: 3/22/2010, 3 hours, QNX 3D implementation design ---
$ws.Range("A62").Copy() | Out-Null
$ws.Range("A64").PasteSpecial(-4122) | Out-Null
$ws.Range("A64").Value = "03/22/2010"
$ws.Range("B64").Value = 3
$ws.Range("C64").Value = "QNX 3D implementation design"

$excel.CutCopyMode = $false

# --- Sheet view updates: scroll to row 37, select A65 ---
$ws.Range("A65").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 37
